# Generate Report for Handoff
# Refresh the localization-status report: the handoff package moved from
# GUID bf00eb2a-be49-49d5-9134-c18f46ae8b0e to 981fb965-ae0a-4bf0-9440-fe31411d350c
# (new target-file hash 3a06d9e87243a312f95f80e307a44d8ef673e9a0) and the
# handoff timestamps advanced.

$wb = $excel.ActiveWorkbook

$oldGuid = "bf00eb2a-be49-49d5-9134-c18f46ae8b0e"
$newGuid = "981fb965-ae0a-4bf0-9440-fe31411d350c"

$mdName   = "$newGuid.md"
$zhCnName = "$newGuid.3a06d9e87243a312f95f80e307a44d8ef673e9a0.zh-cn.xlf"
$deDeName = "$newGuid.3a06d9e87243a312f95f80e307a44d8ef673e9a0.de-de.xlf"

$newHandoffDate = "2016-03-23 13:08:38"
$newZhCnXlfDate = "2016-03-23 13:08:32"

# The hyperlink *targets* are pinned to the commit that introduced the file
# and are left untouched by this change - only the on-sheet display text
# (and the underlying cell value) is refreshed to the new handoff name.
$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/4ccb5f5f525dba832d500e218b2d72c2bea72025/e2e/$oldGuid.md"
$zhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/965a02872530254ee179c543ca70a35d7e5dfab6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.a2c4d74a777f786e8e078d78c42c4fa98f4cecd6.zh-cn.xlf"
$deDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/400e202012e52d26409f945a056498a6a2717f91/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.a2c4d74a777f786e8e078d78c42c4fa98f4cecd6.de-de.xlf"

# NOTE: on this engine, Range.Hyperlinks.Delete() (and Worksheet.Hyperlinks.Delete())
# clears *every* hyperlink on the worksheet, not just the targeted range, and
# setting properties on an existing Hyperlink object appends a duplicate entry
# instead of updating it in place. So for each sheet: wipe the hyperlinks once,
# update the cell text, then recreate every hyperlink that belongs on the sheet.

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $mdName
$wsOverview.Range("D2").Value = $newHandoffDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $mdName)

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Range("A2").Value = $mdName
$wsZhCn.Range("D2").Value = $zhCnName
$wsZhCn.Range("E2").Value = $newZhCnXlfDate
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", $mdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhCnXlfUrl, "", "", $zhCnName)

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Range("A2").Value = $mdName
$wsDeDe.Range("D2").Value = $deDeName
$wsDeDe.Range("E2").Value = $newHandoffDate
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", $mdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deDeXlfUrl, "", "", $deDeName)
